# Added Meeting #4 Minutes / Meeting #4 Proof — Progress&Estimations.xlsx
#
# Adds a 4th meeting ("Mtg 4") column to the "Meeting Attendance" table
# (everyone present), extends the "Hours Worked Per Day" table with an
# extra day's worth of hours (columns G/H), and extends the
# "Task Progress Daily Estimate" figures with one more day (column AF),
# including the QA-check percentage on row 11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Hours Worked Per Day (rows 6-11) -----------------------------------
# Row 6 (Shan): new day columns G/H, plus extra daily progress estimate AF
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 0.5
$ws.Range("AF6").Value = 80

# Row 7 (Pablo)
$ws.Range("H7").Value = 0.5
$ws.Range("AF7").Value = 80

# Row 8 (Jun)
$ws.Range("H8").Value = 0.5
$ws.Range("AF8").Value = 0

# Row 9 (Pedro)
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("AF9").Value = 0

# Row 10 (Brian)
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("AF10").Value = 0

# Row 11 (Young) + QA check % estimate for the new day
$ws.Range("G11").Value = 0.5
$ws.Range("H11").Value = 1.5
$ws.Range("AF11").Value = 0.2

# --- Meeting Attendance table (rows 22-30) ------------------------------
# Header: add the 4th meeting label
$ws.Range("H22").Value = "Mtg 4"

# Attendance for Mtg 4 - everyone present
$ws.Range("H25").Value = "Present"
$ws.Range("H26").Value = "Present"
$ws.Range("H27").Value = "Present"
$ws.Range("H28").Value = "Present"
$ws.Range("H29").Value = "Present"
$ws.Range("H30").Value = "Present"

# --- View state: move selection to where the edits were made -----------
$ws.Range("AG27").Select()
